# Daily attendance processing - 2025-11-15 12:36:08
# Normalizes the "Recorded By" column (G): when a cell lists multiple
# recorders separated by ", ", the first-listed recorder is moved to the
# end of the list (rotate left by one) -- unless the list already starts
# with the exact token "System", in which case the cell is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -ge 2 -and $parts[0] -ne "System") {
        $first = $parts[0]
        $rest = $parts[1..($parts.Count - 1)]
        $newParts = $rest + $first
        $newVal = $newParts -join ", "
        $cell.Value = $newVal
    }
}
